$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing J5 value (date)
$ws.Range("J5").Value = [DateTime]"2024-05-28"

# Copy row 6 formatting down into row 7 (new row) so it reuses existing styles
$ws.Range("B6:J6").Copy()
$ws.Range("B7:J7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B7").Value = [DateTime]"2024-05-28"
$ws.Range("C7").Value = "Configuraciones"
$ws.Range("D7").Value = "Ampliar la idea de configuraciones del sistema donde en vez de solo tener una opcion que sea cambiar el ruteo de conexion se tengan otras opciones adicionales"
$ws.Range("E7").Value = "Medio"
$ws.Range("F7").Value = "Darle al usuario mas libertad de personalizacion del sistema"
$ws.Range("G7").Value = "Se tendra que ampliar la ventana de configuraciones creando las mismas y las funcionalidades correspondientes a cada opcion"
$ws.Range("H7").Value = "-"
$ws.Range("I7").Value = [DateTime]"2024-06-01"
$ws.Rows.Item(7).RowHeight = 72

# Delete old rows 8,9,10 (blank placeholder rows) - shift cells up
$ws.Range("B8:J10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

$ws.Range("L5").Select()
